# Apply "minor template changes, more checkboxes and ui update"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder K1:M1 header values:
#   K1: Herdlaesion -> Lateralitaet
#   L1: Zweitlaesion -> Herdlaesion
#   M1: Lateralitaet -> Zweitlaesion
$ws.Range("K1").Value = "Lateralität"
$ws.Range("L1").Value = "Herdläsion"
$ws.Range("M1").Value = "Zweitläsion"

# Widen the new checkbox column (L, col 12) to match column K (col 11) width
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# Update the view: scroll back into view and select the full column N
$ws.Columns("N:N").Select()
